$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.872.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.49%  "

# Row 3: 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.94%  "

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = "  +0.24%  "

# Row 5: 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6: 'XRP'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.603"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.78%  "

# Row 7: 'USDC'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8: 'Solana'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -8.69%  "

# Row 9: 'Cardano'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.321"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.28%  "

# Row 10: 'Dogecoin'
$ws.Range("E10").Value = "  -2.88%  "

# Row 11: 'TRON'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0992"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.38%  "

# Row 12: 'WrappedliquidstakedEther2.0'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.064.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.99%  "

# Row 13: 'WrappedEther'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.800.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.24%  "

# Row 14: 'Polygon'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.656"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.64%  "

# Row 15: 'Chainlink'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.07%  "

# Row 16: 'Polkadot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.43%  "

# Row 17: 'WrappedBTC'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.788.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.56%  "

# Row 18: 'Litecoin'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.91%  "

# Row 19: 'ShibaInu'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0777"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.35%  "

# Row 20: 'BitcoinCash'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.99%  "

# Row 21: 'Avalanche'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.43%  "

# Row 22: 'Uniswap'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.01%  "

# Row 23: 'Dai'
$ws.Range("E23").Value = "  +0.21%  "

# Row 24: 'Toncoin'
$ws.Range("E24").Value = "  -0.56%  "

# Row 25: 'Monero'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "

# Row 26: 'Cosmos'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.72%  "

# Row 27: 'EthereumClassic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.98%  "

# Row 28: 'Stellar'
$ws.Range("E28").Value = "  -3.65%  "

# Row 29: 'PancakeSwap'
$ws.Range("E29").Value = "  +5.60%  "

# Row 30: 'BinanceUSD'
$ws.Range("E30").Value = "  +0.16%  "

# Row 31: 'Filecoin'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "

# Row 32: 'Hedera'
$ws.Range("E32").Value = "  +0.81%  "

# Row 33: 'InternetComputer(DFINITY)'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.55%  "

# Row 34: 'LidoDAOToken'
$ws.Range("E34").Value = "  -7.96%  "

# Row 35: 'TrustWalletToken'
$ws.Range("E35").Value = "  +4.03%  "

# Row 36: 'ImmutableX'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.678"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.19%  "

# Row 37: 'Aave'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "90.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.39%  "

# Row 38: 'WEMIXToken'
$ws.Range("E38").Value = "  +1.36%  "

# Row 39: 'Maker'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.305.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.80%  "

# Row 40: 'VeChain'
$ws.Range("E40").Value = "  -3.00%  "

# Row 41: 'HuobiToken'
$ws.Range("E41").Value = "  -1.12%  "

# Row 42: 'ARBITRUM' -> 'InjectiveProtocol'
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.89%  "

# Row 43: 'InjectiveProtocol' -> 'ARBITRUM'
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.954"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.18%  "

# Row 44: 'MXToken'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.47%  "

# Row 45: 'RenderToken'
$ws.Range("E45").Value = "  -13.14%  "

# Row 46: 'FraxShare'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "

# Row 47: 'Kaspa'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0508"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.64%  "

# Row 48: 'RocketPoolETH'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.984.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.03%  "

# Row 49: 'PaxDollar'
$ws.Range("E49").Value = "  +0.20%  "

# Row 50: 'Cronos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0669"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.46%  "

# Row 51: 'Quant'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.73%  "
